$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 388.46155
$ws.Range("I38").Value = 295.45456
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 886.36368
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = -514.36368
$ws.Range("N38").Value = -3444
$ws.Range("H88").Value = 4535.4287
$ws.Range("I88").Value = 9925
$ws.Range("J88").Value = 2379.6
$ws.Range("K88").Value = 9925
$ws.Range("L88").Value = 2379.6
$ws.Range("M88").Value = -9519
$ws.Range("N88").Value = -3191.6
$ws.Range("H91").Value = 4535.4287
$ws.Range("I91").Value = 9925
$ws.Range("J91").Value = 2379.6
$ws.Range("K91").Value = 9925
$ws.Range("L91").Value = 2379.6
$ws.Range("M91").Value = -8521
$ws.Range("N91").Value = -5187.6
$ws.Range("H112").Value = 1646.25
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1679
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 5037
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -7253
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442
$ws.Range("H129").Value = 1295.0646
$ws.Range("J129").Value = 1334.4828
$ws.Range("L129").Value = 4003.4484
$ws.Range("N129").Value = -14003.4484
$ws.Range("H130").Value = 60780
$ws.Range("J130").Value = 60780
$ws.Range("L130").Value = 60780
$ws.Range("N130").Value = -70820
$ws.Range("H141").Value = 1975.6666
$ws.Range("I141").Value = 1597.625
$ws.Range("K141").Value = 4792.875
$ws.Range("M141").Value = 387.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 9864.4375
$ws.Range("I32").Value = 8595.071
$ws.Range("J32").Value = 18750
$ws.Range("K32").Value = 8595.071
$ws.Range("L32").Value = 18750
$ws.Range("M32").Value = -8308.071
$ws.Range("N32").Value = -19324
$ws.Range("H37").Value = 7403
$ws.Range("J37").Value = 8353.75
$ws.Range("L37").Value = 8353.75
$ws.Range("N37").Value = -8899.75
$ws.Range("H44").Value = 28833.166
$ws.Range("J44").Value = 28833.166
$ws.Range("L44").Value = 28833.166
$ws.Range("N44").Value = -29809.166
$ws.Range("H55").Value = 30125.309
$ws.Range("J55").Value = 30125.309
$ws.Range("L55").Value = 30125.309
$ws.Range("N55").Value = -30755.309
$ws.Range("H80").Value = 24072.334
$ws.Range("J80").Value = 24072.334
$ws.Range("L80").Value = 24072.334
$ws.Range("N80").Value = -26068.334
$ws.Range("H83").Value = 24072.334
$ws.Range("J83").Value = 24072.334
$ws.Range("L83").Value = 72217.00199999999
$ws.Range("N83").Value = -82201.00199999999
$ws.Range("H132").Value = 4168248.5
$ws.Range("I132").Value = 5001263.5
$ws.Range("J132").Value = 3172.4
$ws.Range("K132").Value = 15003790.5
$ws.Range("L132").Value = 9517.200000000001
$ws.Range("M132").Value = -15001260.5
$ws.Range("N132").Value = -14577.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5006001.5
$ws.Range("I6").Value = 5631501.5
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 5631501.5
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -5631388.5
$ws.Range("N6").Value = -2226
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 7579625
$ws.Range("I31").Value = 4051.8057
$ws.Range("J31").Value = 41669704
$ws.Range("K31").Value = 4051.8057
$ws.Range("L31").Value = 41669704
$ws.Range("M31").Value = -3756.8057
$ws.Range("N31").Value = -41670294
$ws.Range("H34").Value = 7579625
$ws.Range("I34").Value = 4051.8057
$ws.Range("J34").Value = 41669704
$ws.Range("K34").Value = 4051.8057
$ws.Range("L34").Value = 41669704
$ws.Range("M34").Value = -3849.8057
$ws.Range("N34").Value = -41670108
$ws.Range("H41").Value = 12166.667
$ws.Range("I41").Value = 6500
$ws.Range("K41").Value = 6500
$ws.Range("M41").Value = -6072
$ws.Range("H50").Value = 8999
$ws.Range("J50").Value = 8999
$ws.Range("L50").Value = 8999
$ws.Range("N50").Value = -10249
$ws.Range("H59").Value = 16265
$ws.Range("J59").Value = 16265
$ws.Range("L59").Value = 16265
$ws.Range("N59").Value = -18555
$ws.Range("H60").Value = 7520.6665
$ws.Range("I60").Value = 3000
$ws.Range("J60").Value = 8424.799999999999
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 8424.799999999999
$ws.Range("M60").Value = -2489
$ws.Range("N60").Value = -9446.799999999999
$ws.Range("H68").Value = 23725
$ws.Range("J68").Value = 23725
$ws.Range("L68").Value = 23725
$ws.Range("N68").Value = -25223
$ws.Range("H71").Value = 23725
$ws.Range("J71").Value = 23725
$ws.Range("L71").Value = 71175
$ws.Range("N71").Value = -78663
$ws.Range("H74").Value = 13698.6
$ws.Range("J74").Value = 13698.6
$ws.Range("L74").Value = 13698.6
$ws.Range("N74").Value = -15446.6
$ws.Range("H77").Value = 13698.6
$ws.Range("J77").Value = 13698.6
$ws.Range("L77").Value = 41095.8
$ws.Range("N77").Value = -49831.8
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 19232462
$ws.Range("I132").Value = 25001096
$ws.Range("J132").Value = 3685
$ws.Range("K132").Value = 75003288
$ws.Range("L132").Value = 11055
$ws.Range("M132").Value = -75000758
$ws.Range("N132").Value = -16115

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1229.1428
$ws.Range("J25").Value = 1600.8
$ws.Range("L25").Value = 4802.4
$ws.Range("N25").Value = -5140.4
$ws.Range("H30").Value = 1229.1428
$ws.Range("J30").Value = 1600.8
$ws.Range("L30").Value = 4802.4
$ws.Range("N30").Value = -5006.4
$ws.Range("H50").Value = 152.3158
$ws.Range("I50").Value = 127.85714
$ws.Range("J50").Value = 220.8
$ws.Range("K50").Value = 383.57142
$ws.Range("L50").Value = 662.4000000000001
$ws.Range("M50").Value = 97.42858000000001
$ws.Range("N50").Value = -1624.4
$ws.Range("H53").Value = 152.3158
$ws.Range("I53").Value = 127.85714
$ws.Range("J53").Value = 220.8
$ws.Range("K53").Value = 383.57142
$ws.Range("L53").Value = 662.4000000000001
$ws.Range("M53").Value = 97.42858000000001
$ws.Range("N53").Value = -1624.4
$ws.Range("H69").Value = 1283.3158
$ws.Range("I69").Value = 738.1539
$ws.Range("K69").Value = 2214.4617
$ws.Range("M69").Value = -1403.4617
$ws.Range("H72").Value = 1283.3158
$ws.Range("I72").Value = 738.1539
$ws.Range("K72").Value = 6643.3851
$ws.Range("M72").Value = -2587.3851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 33.166668
$ws.Range("I2").Value = 38.5
$ws.Range("J2").Value = 27.833334
$ws.Range("K2").Value = 38.5
$ws.Range("L2").Value = 27.833334
$ws.Range("M2").Value = 74.5
$ws.Range("N2").Value = -253.833334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 800
$ws.Range("I9").Value = 600
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 1200
$ws.Range("M9").Value = -376
$ws.Range("N9").Value = -1648
$ws.Range("H20").Value = 8000
$ws.Range("J20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("N20").Value = -8452
$ws.Range("H132").Value = 11632954
$ws.Range("I132").Value = 3900.3
$ws.Range("J132").Value = 38469230
$ws.Range("K132").Value = 11700.9
$ws.Range("L132").Value = 115407690
$ws.Range("M132").Value = -9170.900000000001
$ws.Range("N132").Value = -115412750
